$wb = $excel.ActiveWorkbook

# --- 1. Fix the shared-string text on sheet "KPI" cell D3: drop the stray
#        space before the comma ("...Location ,Chewing..." -> "...Location,Chewing...")
$kpi = $wb.Worksheets.Item("KPI")
$kpi.Range("D3").Value = "Chewing Gum Primary Location,Chewing Gum Secondary Location"

# --- 2. Workbook-level tab ratio tweak (992 -> 993)
$excel.ActiveWindow.TabRatio = 993

# --- 3. Column width touch-ups (template re-layout to match DB change)
$kpi.Columns.Item(1).ColumnWidth = 31.838731443994565
$kpi.Columns.Item(2).ColumnWidth = 12.235492577597867
$kpi.Columns.Item(3).ColumnWidth = 40.729419703103964
$kpi.Columns.Item(4).ColumnWidth = 56.688933873144364
$kpi.Columns.Item(5).ColumnWidth = 17.057354925775968

$visible = $wb.Worksheets.Item("Visible")
$visible.Columns.Item(2).ColumnWidth = 52.19095816464237

$setSize = $wb.Worksheets.Item("set size")
$setSize.Columns.Item(1).ColumnWidth = 54.44197031039137

# --- 4. Selection clean-up per sheet (drop the stray "D4" from the sqref,
#        collapsing the multi-area selection back down to a single cell)
$visible.Activate() | Out-Null
$visible.Range("B8").Select() | Out-Null

$setSize.Activate() | Out-Null
$setSize.Range("D21").Select() | Out-Null

# Re-activate KPI last so it stays the selected/visible tab, matching the
# unchanged tabSelected="true" on the KPI sheet.
$kpi.Activate() | Out-Null
$kpi.Range("D3").Select() | Out-Null
